$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion-of-the-day message text ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.22 = 8446.67 pesos`n✅ 8446.67 pesos = 2.21 = 947.56 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 450
$ws2.Range("O10").Value = 3801
$ws2.Range("N12").Value = 3821.99
$ws2.Range("O12").Value = 428.755
